$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each cell is written with a leading apostrophe to force Excel to store the
# value as literal text (matching the source data, which are text-formatted
# numbers/percentages, not real numbers), then the cell style is reset to
# "Normal" so the quote-prefix text marker does not leave a stray style on
# the cell (keeping cell formatting identical to the original).

$ws.Range("D2").Value = "'45.384.07"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  +7.20%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'2.385.03"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  +4.83%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  -0.21%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'111.98"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  +8.96%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'317.58"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  +2.82%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("E7").Value = "'  +2.89%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("E8").Value = "'  -0.16%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("E9").Value = "'  +5.64%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'42.36"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  +10.54%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("E11").Value = "'  +4.20%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("E12").Value = "'  +6.41%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("E13").Value = "'  +5.75%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("E14").Value = "'  +1.43%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'15.87"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  +6.36%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'2.745.93"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  +4.85%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'2.385.18"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  +4.89%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'45.343.78"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  +7.20%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("E19").Value = "'  +6.94%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("E20").Value = "'  +4.71%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'13.17"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  +2.23%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'75.06"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  +3.59%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("E23").Value = "'  +5.56%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'269.71"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  +3.31%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'2.33"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  +8.43%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("E26").Value = "'  -0.75%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("E27").Value = "'  +6.85%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'7.53"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  +10.93%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("E29").Value = "'  +0.33%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'22.92"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  +4.14%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'38.55"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  +8.75%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'0.0944"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  +11.46%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'169.81"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  +3.91%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("E34").Value = "'  +17.71%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("E35").Value = "'  +3.63%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("E36").Value = "'  +7.00%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "'4.89"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  +9.90%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = "'3.12"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  +14.90%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("E39").Value = "'  +6.13%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("E40").Value = "'  +9.00%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("E41").Value = "'  +13.37%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'105.27"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  +6.85%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'0.240"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  +7.50%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'13.55"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  +14.93%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'71.17"
$ws.Range("D45").Style = "Normal"
$ws.Range("E46").Value = "'  +0.18%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'117.83"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  +8.15%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'5.81"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  +14.32%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'9.38"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  +9.61%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'1.63"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  +20.23%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'79.06"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  +5.08%  "
$ws.Range("E51").Style = "Normal"
